$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Foreign Currencies")
$summary = $wb.Worksheets.Item("ELSTER - Summary")

# --- Fix up the date-label cells first (while the donor cells still hold
#     their original values) by copying existing shared-string cells rather
#     than re-typing the dates, so Excel doesn't auto-convert the text into
#     a date serial number / introduce new cell styles ---
$ws.Range("C5").Copy($ws.Range("C6"))   # "2022-09-22" -> row 6
$ws.Range("D7").Copy($ws.Range("D6"))   # "2022-12-01" -> row 6
$ws.Range("D7").Copy($ws.Range("D5"))   # "2022-12-01" -> row 5
$ws.Range("C4").Copy($ws.Range("C5"))   # "2022-09-05" -> row 5

# --- Update row 4 (B4, G4) ---
$ws.Range("B4").Value = 2582.03
$ws.Range("G4").Value = 20.72

# --- Update row 5 (B5, F5, G5) ---
$ws.Range("B5").Value = 849.9400000000001
$ws.Range("F5").Value = 0.9399999999999999
$ws.Range("G5").Value = -55.52

# --- Update row 6 (B6, F6, G6) ---
$ws.Range("B6").Value = 150.06
$ws.Range("F6").Value = 0.9399999999999999
$ws.Range("G6").Value = -9.800000000000001

# --- Remove the now-obsolete individual withdrawal rows (old rows 7-11); the
#     summary rows below them shift up to become rows 7-10 ---
$ws.Rows("7:11").Delete()

# --- Update the totals (now rows 8-10) to match the recomputed withdrawal
#     totals above ---
$ws.Range("G8").Value = -33.66
$ws.Range("G9").Value = 31.66
$ws.Range("G10").Value = -65.31999999999999

# --- Mirror the "Gains (incl. losses)" total into the ELSTER summary sheet ---
$summary.Range("C7").Value = -33.66
